# Iterationsplan_as223jx.xlsx - "Lagt till Iterationsplan för Iteration 2"
#
# Summary of the change:
#  - Sheet "Iteration 1": the task "Fortsätta arbeta med testapplikationen" is
#    replaced by "Fixa textruta och action bar på testapplikationen", several
#    "Verklig tid" (actual time) values get filled in, and the running totals
#    are updated.
#  - Sheet3 gets a brand-new "Iterationsplan Iteration 2" content block (mirrors
#    the layout already used by "Iteration 0" / "Iteration 1"), and becomes the
#    active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. "Iteration 1" sheet updates
# ---------------------------------------------------------------------------

# The task text in row 14 changed (this also happens to be the first newly
# introduced shared string, so it is written before anything on Sheet3).
$ws2.Range("A14").Value = "Fixa textruta och action bar på testapplikationen"

# Newly-filled "Verklig tid" (actual time) values.
$ws2.Range("D13").Value = 0
$ws2.Range("D14").Value = 2
$ws2.Range("D17").Value = 1
$ws2.Range("D18").Value = 1

# Updated sums.
$ws2.Range("D20").Value = 13
$ws2.Range("D21").Value = 216

# Column A grew a bit wider to fit the new text.
$ws2.Columns.Item(1).ColumnWidth = 44.7109375

# ---------------------------------------------------------------------------
# 2. Build the new "Sheet3" content (Iterationsplan Iteration 2)
# ---------------------------------------------------------------------------

$ws3.Range("A1").Value = "Iterationsplan Iteration 2"
$ws3.Range("A3").Value = "Jobbat med prototyp och persona, samt börjat implementera textfält, knappar och action bar i min testapplikation"
$ws3.Range("A5").Value = "Målet med denna iteration är att börja implementera grundläggande funktioner i min egen applikation."

$ws3.Range("A13").Value = "Skapa sökruta"
$ws3.Range("A9").Value  = "Skapa nytt projekt i Eclipse för min app"
$ws3.Range("A10").Value = "Skapa ""Nytt recept""-knapp som öppnar ny sida"
$ws3.Range("A14").Value = "Planera färgtema & bakgrund"
$ws3.Range("A15").Value = "Designa ikon till appen"
$ws3.Range("A11").Value = "Skapa textfälten samt sparaknapp i ""Nytt recept""-sidan"
$ws3.Range("A16").Value = "Skapa Iterationsplan för Iteration 3"

# Cells that reuse strings already present elsewhere in the workbook.
$ws3.Range("A2").Value = "Analys av föregående iteration"
$ws3.Range("A4").Value = "Mål"

$ws3.Range("A7").Value = "Uppgift"
$ws3.Range("B7").Value = "Status"
$ws3.Range("C7").Value = "Skattad tid"
$ws3.Range("D7").Value = "Verklig tid"
$ws3.Range("E7").Value = "Kommentar"

$ws3.Range("A8").Value  = "Handledarmöte"
$ws3.Range("A12").Value = "Finslipning av krav"

$ws3.Range("B17").Value = "Summa"
$ws3.Range("B18").Value = "Tid sedan föregående iteration"
$ws3.Range("B19").Value = "Total projekttid"

# Status column ("Ej påbörjad" for every task row) + estimated-time column.
$taskRows = 8..16
foreach ($r in $taskRows) {
    $ws3.Range("B$r").Value = "Ej påbörjad"
}
$ws3.Range("C8").Value  = 1
$ws3.Range("C9").Value  = 1
$ws3.Range("C10").Value = 1
$ws3.Range("C11").Value = 3
$ws3.Range("C12").Value = 1
$ws3.Range("C13").Value = 1
$ws3.Range("C14").Value = 2
$ws3.Range("C15").Value = 1
$ws3.Range("C16").Value = 1

$ws3.Range("C17").Value = 12
$ws3.Range("D18").Value = 216
$ws3.Range("D19").Value = 240

# ---------------------------------------------------------------------------
# 3. Formatting: reuse the same cell styles the other iteration sheets use.
# ---------------------------------------------------------------------------

# Title row (dark-grey, centered) - same style as "Iteration 1" row 1.
$ws2.Range("A1:E1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)
$ws3.Range("A1:E1").Merge()

# "Analys av föregående iteration" section header (light-grey, centered).
$ws2.Range("A2:E2").Copy()
$ws3.Range("A2:E2").PasteSpecial(-4122)
$ws3.Range("A2:E2").Merge()

# "Mål" section header (light-grey, centered).
$ws2.Range("A4:E4").Copy()
$ws3.Range("A4:E4").PasteSpecial(-4122)
$ws3.Range("A4:E4").Merge()

# Goal text box (wrap text only, spans A:D).
$ws2.Range("A5:D5").Copy()
$ws3.Range("A5:D5").PasteSpecial(-4122)
$ws3.Range("A5:D5").Merge()

# Analysis/comment text box - wrap text, but also centered (new style),
# and it spans the full A:E width on this sheet.
$ws2.Range("A3:D3").Copy()
$ws3.Range("A3:E3").PasteSpecial(-4122)
$ws3.Range("A3:E3").HorizontalAlignment = -4108
$ws3.Range("A3:E3").Merge()

# Table header row (bold).
$ws2.Range("A7:E7").Copy()
$ws3.Range("A7:E7").PasteSpecial(-4122)

# Summary rows (grey label cells).
$ws2.Range("B20").Copy()
$ws3.Range("B17").PasteSpecial(-4122)
$ws2.Range("B21").Copy()
$ws3.Range("B18").PasteSpecial(-4122)
$ws2.Range("B22").Copy()
$ws3.Range("B19").PasteSpecial(-4122)

$ws3.Columns.Item(1).ColumnWidth = 50.7109375

# ---------------------------------------------------------------------------
# 4. Selections / active sheet
# ---------------------------------------------------------------------------

$ws2.Range("D21").Select()
$ws3.Activate()
$ws3.Range("D20").Select()
